# Update "Forecast Comparison" sheet:
#  - insert a new "Week_Start_Date" column between "Week" and "ASIN"
#  - normalize the Week labels from W01..W16 to W1..W16
#  - populate the new column with the week start dates (as text)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new blank column before column B (ASIN etc. shift right by one).
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Week start dates (Sundays) for weeks 1-16 of 2025, aligned to rows 2-17.
$weekStarts = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

# Make sure the new column stores the dates as plain text (matching the
# source data, which keeps them as literal "YYYY-MM-DD" strings rather
# than Excel date serials).
$ws.Range("B2:B17").NumberFormat = "@"

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2

    # Normalize "W01".."W16" -> "W1".."W16"
    $weekNum = $i + 1
    $ws.Cells.Item($row, 1).Value = "W" + $weekNum

    # Fill in the new Week_Start_Date column.
    $ws.Cells.Item($row, 2).Value = $weekStarts[$i]
}
